# pest control file object now support observation groups and parameter groups
#
# Updates the "PST file" sheet's "usages" column (B) for the rows that track
# whether the parser/writer object supports a given PEST control-file
# section. "parameter groups" and "parameter data" go from unsupported (0)
# to fully supported (1) and gain a "BeoJACTEST" note in column C (matching
# the note already present on "control data"); "observation groups" goes
# from 0 to 1 as well. The two rows that were already partially supported
# ("block separation", "control data") become fully supported (1).
#
# The ParserWriter sheet's summary formula in C6 (AVERAGE of the "PST file"
# usages column) recalculates automatically as a result.

$wb = $excel.ActiveWorkbook
$pstSheet = $wb.Worksheets.Item("PST file")
$pwSheet  = $wb.Worksheets.Item("ParserWriter")

# block separation: 0.5 -> 1
$pstSheet.Range("B5").Value = 1

# control data: 0.5 -> 1
$pstSheet.Range("B7").Value = 1

# parameter groups: 0 -> 1, now implemented via BeoJACTEST
$pstSheet.Range("B13").Value = 1
$pstSheet.Range("C13").Value = "BeoJACTEST"

# parameter data: 0 -> 1, now implemented via BeoJACTEST
$pstSheet.Range("B14").Value = 1
$pstSheet.Range("C14").Value = "BeoJACTEST"

# observation groups: 0 -> 1
$pstSheet.Range("B15").Value = 1

# Move the active selection/tab so the workbook reopens focused on the
# ParserWriter sheet at C6, while the PST file sheet keeps a parked
# selection at B16.
$pstSheet.Activate()
$pstSheet.Range("B16").Select()

$pwSheet.Activate()
$pwSheet.Range("C6").Select()
